# Apply weekly update of Fruta/Hortaliza (Caqui) data.
# The edit re-shuffles the Date/Calidad/Volumen/Precio/Unidad/Origen/Precio-Kg/Kg-unidad
# values across rows 2, 3, 5, 7 and 8 (rows 4 and 6 are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2, 4).Value  = 44714                          # D2 Fecha
$ws.Cells.Item(2, 14).Value = 14000                           # N2 Precio minimo
$ws.Cells.Item(2, 15).Value = 15000                           # O2 Precio maximo
$ws.Cells.Item(2, 16).Value = 14500                           # P2 Precio promedio ponderado
$ws.Cells.Item(2, 18).Value = "Región de O'Higgins"           # R2 Origen
$ws.Cells.Item(2, 19).Value = 806                             # S2 Precio $/Kg

# --- Row 3 ---
$ws.Cells.Item(3, 4).Value  = 44707                           # D3 Fecha
$ws.Cells.Item(3, 12).Value = "Primera"                       # L3 Calidad
$ws.Cells.Item(3, 13).Value = 60                              # M3 Volumen
$ws.Cells.Item(3, 14).Value = 12000                           # N3 Precio minimo
$ws.Cells.Item(3, 15).Value = 13000                           # O3 Precio maximo
$ws.Cells.Item(3, 16).Value = 12500                           # P3 Precio promedio ponderado
$ws.Cells.Item(3, 17).Value = "$/caja 12 kilos empedrada"     # Q3 Unidad de comercializacion
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"           # R3 Origen
$ws.Cells.Item(3, 19).Value = 1042                            # S3 Precio $/Kg
$ws.Cells.Item(3, 20).Value = 12                              # T3 Kg / unidad

# --- Row 5 ---
$ws.Cells.Item(5, 4).Value  = 44742                           # D5 Fecha
$ws.Cells.Item(5, 12).Value = "Segunda"                       # L5 Calidad
$ws.Cells.Item(5, 14).Value = 14000                           # N5 Precio minimo
$ws.Cells.Item(5, 15).Value = 15000                           # O5 Precio maximo
$ws.Cells.Item(5, 16).Value = 14500                           # P5 Precio promedio ponderado
$ws.Cells.Item(5, 17).Value = "$/caja 18 kilos granel"        # Q5 Unidad de comercializacion
$ws.Cells.Item(5, 19).Value = 806                             # S5 Precio $/Kg
$ws.Cells.Item(5, 20).Value = 18                              # T5 Kg / unidad

# --- Row 7 ---
$ws.Cells.Item(7, 4).Value  = 44330                           # D7 Fecha
$ws.Cells.Item(7, 13).Value = 100                             # M7 Volumen
$ws.Cells.Item(7, 14).Value = 15000                           # N7 Precio minimo
$ws.Cells.Item(7, 15).Value = 16000                           # O7 Precio maximo
$ws.Cells.Item(7, 16).Value = 15500                           # P7 Precio promedio ponderado
$ws.Cells.Item(7, 17).Value = "$/caja 18 kilos granel"        # Q7 Unidad de comercializacion
$ws.Cells.Item(7, 19).Value = 861                             # S7 Precio $/Kg
$ws.Cells.Item(7, 20).Value = 18                              # T7 Kg / unidad

# --- Row 8 ---
$ws.Cells.Item(8, 4).Value  = 44334                           # D8 Fecha
$ws.Cells.Item(8, 14).Value = 11000                           # N8 Precio minimo
$ws.Cells.Item(8, 15).Value = 12000                           # O8 Precio maximo
$ws.Cells.Item(8, 16).Value = 11500                           # P8 Precio promedio ponderado
$ws.Cells.Item(8, 17).Value = "$/caja 12 kilos granel"        # Q8 Unidad de comercializacion
$ws.Cells.Item(8, 19).Value = 11500                           # S8 Precio $/Kg
$ws.Cells.Item(8, 20).Value = 1                               # T8 Kg / unidad
